$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column C (copy the header format used by A1/B1)
$ws.Range("C1").Value = 'articletitle'
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Article title values for rows 2-21
$ws.Range("C2").Value = 'Books that explain the world: Guardian writers share their best nonfiction reads of the year'
$ws.Range("C3").Value = 'Could he? Will they? What if? What might happen next in Succession'
$ws.Range("C4").Value = 'Rajan the last ocean-swimming elephant: Jody MacDonald’s best photograph'
$ws.Range("C5").Value = '‘I’ve healed. I don’t want to be the badass’ – Noomi Rapace on beating her Dragon Tattoo trauma'
$ws.Range("C6").Value = 'Shatner in Space: Captain Kirk’s in a willy-shaped spaceship – and it’s poetry in motion'
$ws.Range("C7").Value = 'I lost more than £12,000 in a scam and the Co-operative Bank won’t resolve it'
$ws.Range("C8").Value = 'Are we witnessing the dawn of post-theory science?'
$ws.Range("C9").Value = '‘A Rosetta Stone’: Australian fossil site is a vivid window into 15m-year-old rainforest'
$ws.Range("C10").Value = 'I''m still repaying my repaid student loan'
$ws.Range("C11").Value = 'England’s new attack stutters as West Indies’ Holder and Bonner stand firm'
$ws.Range("C12").Value = 'Spain reverses plan to open up to unvaccinated British visitors'
$ws.Range("C13").Value = 'All the flights cancelled today from UK airports'
$ws.Range("C14").Value = 'Government defeated over voter ID plans in House of Lords'
$ws.Range("C15").Value = '‘We’ve been trying warn you for so many decades’: Nasa climate scientist breaks down in tears at protest'
$ws.Range("C16").Value = 'Student who went to A&E with stomach cramps gives birth to healthy baby boy'
$ws.Range("C17").Value = 'Theresa May: We''re on course to deliver Brexit despite vote'
$ws.Range("C18").Value = 'This season at Nats Park: Plenty of good seats and a quest to find joy'
$ws.Range("C19").Value = 'Cristiano Ronaldo scores landmark goal days after death of infant son'
$ws.Range("C20").Value = 'Division III basketball player apologizes for brutal sucker punch'
$ws.Range("C21").Value = 'Ukraine war: World Bank warns of ''human catastrophe'' food crisis'

# Column width for C
$ws.Columns.Item(3).ColumnWidth = 12.1667

# Special bold/large font style for C12:C13 (featured rows)
$ws.Range("C12:C13").Font.Bold = $true
$ws.Range("C12:C13").Font.Size = 24
$ws.Range("C12:C13").VerticalAlignment = -4108

# Update selection / view state
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B27").Select()
